# Updates cryptos list figures (Price / Volume(1h)) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.124.65'
$ws.Range("E2").Value = '  +3.40%  '
$ws.Range("D3").Value = '3.034.84'
$ws.Range("E3").Value = '  +2.09%  '
$ws.Range("D4").Formula = "'0.999"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Formula = "'592.99"
$ws.Range("E5").Value = '  -0.45%  '
$ws.Range("D6").Formula = "'154.34"
$ws.Range("E6").Value = '  +8.62%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").Value = '3.030.11'
$ws.Range("E8").Value = '  +1.97%  '
$ws.Range("D9").Formula = "'0.517"
$ws.Range("D10").Formula = "'6.85"
$ws.Range("E10").Value = '  +14.29%  '
$ws.Range("E11").Value = '  +4.61%  '
$ws.Range("D12").Formula = "'0.463"
$ws.Range("E12").Value = '  +2.20%  '
$ws.Range("E13").Value = '  +3.69%  '
$ws.Range("D14").Formula = "'35.77"
$ws.Range("E14").Value = '  +5.23%  '
$ws.Range("E15").Value = '  +0.56%  '
$ws.Range("D16").Value = '3.534.15'
$ws.Range("E16").Value = '  +2.06%  '
$ws.Range("D17").Formula = "'7.09"
$ws.Range("E17").Value = '  +3.31%  '
$ws.Range("D18").Value = '63.006.55'
$ws.Range("E18").Value = '  +3.27%  '
$ws.Range("D19").Value = '3.032.30'
$ws.Range("E19").Value = '  +2.02%  '
$ws.Range("D20").Formula = "'453.72"
$ws.Range("E20").Value = '  +1.66%  '
$ws.Range("D21").Formula = "'14.28"
$ws.Range("E21").Value = '  +1.44%  '
$ws.Range("D22").Formula = "'0.698"
$ws.Range("E22").Value = '  +3.03%  '
$ws.Range("D23").Formula = "'7.51"
$ws.Range("E23").Value = '  +3.37%  '
$ws.Range("D24").Formula = "'83.15"
$ws.Range("E24").Value = '  +1.52%  '
$ws.Range("D25").Formula = "'11.34"
$ws.Range("E25").Value = '  +10.11%  '
$ws.Range("D26").Formula = "'2.31"
$ws.Range("E26").Value = '  +6.58%  '
$ws.Range("D27").Formula = "'12.46"
$ws.Range("E27").Value = '  +5.02%  '
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("D29").Formula = "'7.51"
$ws.Range("E29").Value = '  +6.08%  '
$ws.Range("D30").Formula = "'2.26"
$ws.Range("E30").Value = '  +11.26%  '
$ws.Range("E31").Value = '  +1.01%  '
$ws.Range("E32").Value = '  -0.07%  '
$ws.Range("D33").Formula = "'27.60"
$ws.Range("E33").Value = '  +2.17%  '
$ws.Range("E34").Value = '  +2.44%  '
$ws.Range("D35").Value = '0.0₃0863'
$ws.Range("E35").Value = '  +6.43%  '
$ws.Range("E36").Value = '  +3.50%  '
$ws.Range("D37").Formula = "'5.93"
$ws.Range("E37").Value = '  +3.45%  '
$ws.Range("D38").Formula = "'3.18"
$ws.Range("E38").Value = '  +12.08%  '
$ws.Range("D39").Formula = "'0.131"
$ws.Range("E39").Value = '  +8.12%  '
$ws.Range("D40").Formula = "'2.10"
$ws.Range("E40").Value = '  +3.06%  '
$ws.Range("D41").Formula = "'50.53"
$ws.Range("E41").Value = '  +0.68%  '
$ws.Range("D42").Formula = "'9.14"
$ws.Range("E42").Value = '  +1.75%  '
$ws.Range("D43").Formula = "'0.308"
$ws.Range("E43").Value = '  +16.19%  '
$ws.Range("D44").Formula = "'43.97"
$ws.Range("E44").Value = '  +12.80%  '
$ws.Range("D45").Formula = "'391.34"
$ws.Range("E45").Value = '  +0.28%  '
$ws.Range("E46").Value = '  +3.86%  '
$ws.Range("D47").Value = '2.722.52'
$ws.Range("E47").Value = '  +1.70%  '
$ws.Range("D48").Formula = "'133.49"
$ws.Range("E48").Value = '  +2.47%  '
$ws.Range("E50").Value = '  +8.14%  '
$ws.Range("D51").Formula = "'25.15"
$ws.Range("E51").Value = '  +8.70%  '
